$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Subashmahat35@gmail.com"
$ws.Range("D9").Select()
